# Bold five heading/title runs that were left with an explicit
# "bold off" run-level override (w:b w:val="0") even though their
# paragraph style / sibling formatting is bold.

$d = $word.ActiveDocument

$targets = @(
    "Holly Dickson",
    "工作经历",
    "初级动画设计师",
    "动画实习生",
    "动画美术学士学位"
)

foreach ($t in $targets) {
    $rng = $d.Content
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
    }
}
